$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance Changes")

# The "Mark" collection row's example query was replaced with the actual
# mongoose calls used for the update (bulkWrite + find).
$ws.Range("C7").Value = "1 Marks.bulkWrite(updates)     2.  Marks.find(match)"

# Column C no longer needs to be as wide now that the cell holds a much
# shorter string - narrow it down (132 -> 98 characters).
# (ColumnWidth adds Excel's standard 5/6-character padding on save, so we
# back that out here to land on an exact stored width of 98.)
$ws.Columns("C").ColumnWidth = 98 - 5/6

# Rows 4, 6 and 7 pick up the sheet's standard 13.2pt row height (matching
# rows 1 and 11) now that row 7 no longer wraps onto multiple lines.
$ws.Rows("4").RowHeight = 13.2
$ws.Rows("6").RowHeight = 13.2
$ws.Rows("7").RowHeight = 13.2

# Leave the cursor on the edited cell, like the author did before saving.
$ws.Range("C7").Select()
